# Add a new column "Rprc_TH" (= %tree-hedge cover relative to dvlp+TH only)
# to the "huc12" worksheet: a header in S1 plus a per-row ratio formula
# G/(G+H)*100 for rows 2-9, mirroring the existing neighbouring % columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("huc12")

$ws.Range("S1").Value = "Rprc_TH"

for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 19).Formula = "=G$row/(G$row+H$row) *100"
}
